$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2119
$wsExhibit.Range("F4").Value = 871
$wsExhibit.Range("F5").Value = 1331

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2119
$wsAll.Range("F6").Value = 871
$wsAll.Range("F7").Value = 1331
